# Daily attendance processing - 2025-10-26 05:44:12
# Swap the order of the "Recorded By" entries in column G from
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
# wherever that exact value is found.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G
    $current = $cell.Value()
    if ($current -eq $oldValue) {
        $cell.Value = $newValue
    }
}
